# Updated cryptos list on Sat Dec 30 11:12:51 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row, and
# reflects the new relative ranking/ordering for a few coins whose Price
# changes moved them past a neighboring coin (Litecoin/BitcoinCash,
# Toncoin/Cosmos, WEMIXToken/Filecoin swapped places while keeping their
# original rank number in column A).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.068.28'
$ws.Range('E2').Value = '  -1.79%  '

$ws.Range('D3').Value = '2.290.76'
$ws.Range('E3').Value = '  -3.36%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.13'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.15%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.63'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.24%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.627'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.05%  '

$ws.Range('E8').Value = '  +0.13%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.602'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.34%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.06'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.72%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0900'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.18%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.26'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.29%  '

$ws.Range('E13').Value = '  -0.13%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.959'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.22%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.20'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.86%  '

$ws.Range('D16').Value = '2.635.40'
$ws.Range('E16').Value = '  -3.31%  '

$ws.Range('D17').Value = '2.283.71'
$ws.Range('E17').Value = '  -4.44%  '

$ws.Range('D18').Value = '41.947.71'
$ws.Range('E18').Value = '  -2.06%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.43'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.81%  '

$ws.Range('E20').Value = '  -1.40%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.65'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.14%  '

$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '282.63'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +8.88%  '

$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '73.27'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.86%  '

$ws.Range('E24').Value = '  +8.15%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.25'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.51%  '

$ws.Range('E26').Value = '  +0.70%  '

$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.76'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.09%  '

$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.34'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.74%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '22.93'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.30%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.55'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.46%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '163.37'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.13%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0872'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.61%  '

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.83'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.64%  '

$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.84'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.83%  '

$ws.Range('E35').Value = '  +0.98%  '

$ws.Range('E36').Value = '  -3.05%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.54'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.75%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.87'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.95%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0347'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.88%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.67'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.30%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '99.46'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +9.81%  '

$ws.Range('E42').Value = '  -4.00%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '69.23'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.61%  '

$ws.Range('E44').Value = '  +0.14%  '

$ws.Range('E45').Value = '  -7.34%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '114.47'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.92%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '11.88'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.30%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '76.20'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.20%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.92'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.36%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.26'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.55%  '

$ws.Range('E51').Value = '  -2.19%  '
